# Fruta / hortaliza, semanal
# A new weekly observation is inserted as row 26 (Vega Modelo de Temuco -
# Maracuyá), pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 26, shifting rows 26:107 down to 27:108.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(26, 1).Value = 10
$ws.Cells.Item(26, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value = "La Araucanía"
$ws.Cells.Item(26, 4).Value = 45133
$ws.Cells.Item(26, 5).Value = 9
$ws.Cells.Item(26, 6).Value = "Fruta"
$ws.Cells.Item(26, 7).Value = 100108
$ws.Cells.Item(26, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(26, 9).Value = 100108003
$ws.Cells.Item(26, 10).Value = "Maracuyá"
$ws.Cells.Item(26, 11).Value = "Sin especificar"
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 20
$ws.Cells.Item(26, 14).Value = 45000
$ws.Cells.Item(26, 15).Value = 45000
$ws.Cells.Item(26, 16).Value = 45000
$ws.Cells.Item(26, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(26, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 19).Value = 2500
$ws.Cells.Item(26, 20).Value = 18
